# "all fall 22 week 6 inputs complete"
# Fill in the Week-of-9/27 (F) and Week-of-10/4 (G) result columns for both
# tables on the "Wookie Mistakes" sheet — these cells previously held the
# placeholder "A" (Available) and now hold the actual game result.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Wookie Mistakes")
$ws.Activate()

# --- Table 1 (rows 3-10) ---
$ws.Range("F3").Value = "W"
$ws.Range("G3").Value = "W"

$ws.Range("F4").Value = "W"
$ws.Range("G4").Value = "W"

$ws.Range("F5").Value = "W"
$ws.Range("G5").Value = "W"

$ws.Range("F6").Value = "NA"
$ws.Range("G6").Value = "W"

$ws.Range("F7").Value = "L"
$ws.Range("G7").Value = "DNP"

$ws.Range("F8").Value = "W"
$ws.Range("G8").Value = "DNP"

$ws.Range("F9").Value = "NA"
$ws.Range("G9").Value = "W"

$ws.Range("F10").Value = "DNP"
$ws.Range("G10").Value = "DNP"

# --- Table 2 (rows 15-22) ---
$ws.Range("F15").Value = "NA"
$ws.Range("G15").Value = "L"

$ws.Range("F16").Value = "L"
$ws.Range("G16").Value = "DNP"

$ws.Range("F17").Value = "L"
$ws.Range("G17").Value = "W"

$ws.Range("F18").Value = "DNP"
$ws.Range("G18").Value = "L"

$ws.Range("F19").Value = "L"
$ws.Range("G19").Value = "W"

$ws.Range("F20").Value = "W"
$ws.Range("G20").Value = "NA"

$ws.Range("F21").Value = "DNP"
$ws.Range("G21").Value = "DNP"

$ws.Range("F22").Value = "W"
$ws.Range("G22").Value = "L"

# Move the active selection to K28, matching the sheet view after entry.
$ws.Range("K28").Select()
